$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- EntityCard header row (row 4): insert a new "deleted" column before
# "physical power", shifting the existing physical/magical power &
# protecction + PreccenseId columns one cell to the right (F4:J4 -> G4:K4).
#
# Copy J4's current formatting (the "special" PreccenseId style) onto the
# new K4 cell before we overwrite J4 with the shifted-in value.
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K4").Value = "PreccenseId"

# J4 takes over I4's plain bordered style (it is no longer the last/
# "PreccenseId" column), then gets the shifted-in value.
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("J4").Value = "magical protecction"

$ws.Range("I4").Value = "physical protecction"
$ws.Range("H4").Value = "magical power"
$ws.Range("G4").Value = "physical power"
$ws.Range("F4").Value = "deleted"

# --- SkillCard row (row 8): rename "power (from potential)" to "deleted"
# in place (no column shift here; PreccenseId/TypeId stay put).
$ws.Range("F8").Value = "deleted"

# --- Update the view to match the author's saved selection/scroll state.
$ws.Range("F8").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
